{"js": "// Load the single table in the document.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Helper: replace the full text of a single-run table cell (identified by\n// row index, column 0) with a new value while preserving the run's\n// formatting (rFonts/sz). We scope the search to the cell body so\n// duplicate values elsewhere in the table are not touched.\nasync function setCellText(rowIndex, oldText, newText) {\n  const cell = table.getCell(rowIndex, 0);\n  const results = cell.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Helper: a row whose single cell actually holds several tab-separated\n// values packed into one run-sequence (<w:t>.../<w:tab/>/<w:t>.../...).\n// Collapse the whole cell down to a single value, keeping the original\n// formatting of the first run (only its text changes) and dropping the\n// rest of the run/tab sequence.\nasync function collapseRowToValue(rowIndex, firstOldText, newText) {\n  const cell = table.getCell(rowIndex, 0);\n  const results = cell.body.search(firstOldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  const firstMatch = results.items[0];\n\n  const afterFirst = firstMatch.getRange(\"After\");\n  const para = cell.body.paragraphs.getFirst();\n  const paraEnd = para.getRange(\"End\");\n  const tail = afterFirst.expandTo(paraEnd);\n  tail.delete();\n  await context.sync();\n\n  firstMatch.insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Simple single-value cell updates (row indices are 0-based).\nawait setCellText(0, \"99.9\", \"0M\");\nawait setCellText(1, \"0.2\", \"0M\");\nawait setCellText(2, \"201\", \"0M\");\nawait setCellText(3, \"607\", \"1803\");\nawait setCellText(4, \"0.00005\", \"0.00001\");\nawait setCellText(5, \"0.00041\", \"0.00070\");\nawait setCellText(6, \"0.00010\", \"0.00011\");\n// row 7 (0.00004) is unchanged\nawait setCellText(8, \"0.00009\", \"0.00015\");\nawait setCellText(9, \"0.00009\", \"0.00015\");\nawait setCellText(10, \"0.00010\", \"0.00017\");\nawait setCellText(11, \"0.06985\", \"0.20419\");\n\n// The last three rows previously packed ten tab-separated values into a\n// single run-sequence; they now hold just one value each.\nawait collapseRowToValue(43, \"598\", \"99.9\");\nawait collapseRowToValue(44, \"589\", \"0.2\");\nawait collapseRowToValue(45, \"9\", \"201\");\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Simple single-value cells (1-based row numbers within the one-column table).\n$t.Cell(1, 1).Range.Text  = \"0M\"\n$t.Cell(2, 1).Range.Text  = \"0M\"\n$t.Cell(3, 1).Range.Text  = \"0M\"\n$t.Cell(4, 1).Range.Text  = \"1803\"\n$t.Cell(5, 1).Range.Text  = \"0.00001\"\n$t.Cell(6, 1).Range.Text  = \"0.00070\"\n$t.Cell(7, 1).Range.Text  = \"0.00011\"\n# Row 8 (0.00004) is unchanged.\n$t.Cell(9, 1).Range.Text  = \"0.00015\"\n$t.Cell(10, 1).Range.Text = \"0.00015\"\n$t.Cell(11, 1).Range.Text = \"0.00017\"\n$t.Cell(12, 1).Range.Text = \"0.20419\"\n\n# The last three rows previously packed ten tab-separated values into one\n# cell; replacing Range.Text collapses them down to the single new value\n# while keeping the existing run formatting (rFonts/sz).\n$t.Cell(44, 1).Range.Text = \"99.9\"\n$t.Cell(45, 1).Range.Text = \"0.2\"\n$t.Cell(46, 1).Range.Text = \"201\"\n"}
